# Weekly update: insert a new "Coco" price record as the new row 46
# (Mercado Mayorista Lo Valledor de Santiago), pushing the existing
# rows 46-56 down to 47-57.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 46, shifting rows 46:56 down to 47:57.
$ws.Rows.Item(46).Insert()

# Populate the new row 46 with the latest weekly record.
$ws.Cells.Item(46, 1).Value2 = 6
$ws.Cells.Item(46, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(46, 3).Value2 = "Metropolitana"
$ws.Cells.Item(46, 4).Value2 = 44466
$ws.Cells.Item(46, 5).Value2 = 13
$ws.Cells.Item(46, 6).Value2 = "Fruta"
$ws.Cells.Item(46, 7).Value2 = 100108
$ws.Cells.Item(46, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(46, 9).Value2 = 100108007
$ws.Cells.Item(46, 10).Value2 = "Coco"
$ws.Cells.Item(46, 11).Value2 = "Sin especificar"
$ws.Cells.Item(46, 12).Value2 = "Primera"
$ws.Cells.Item(46, 13).Value2 = 70
$ws.Cells.Item(46, 14).Value2 = 19000
$ws.Cells.Item(46, 15).Value2 = 20000
$ws.Cells.Item(46, 16).Value2 = 19500
$ws.Cells.Item(46, 17).Value2 = "$/malla 20 unidades"
$ws.Cells.Item(46, 18).Value2 = "Perú"
$ws.Cells.Item(46, 19).Value2 = 975
$ws.Cells.Item(46, 20).Value2 = 20

# Keep the date column's existing number format (inherited from the row
# insert), and make sure the dimension/used-range reflects the new row 57.
$ws.Cells.Item(46, 4).NumberFormat = $ws.Cells.Item(47, 4).NumberFormat
